$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District column (G) for Uttara Kannada rows (3-21)
$ws.Range("G3:G21").Value = "Uttara Kannada (Karwar)"

# Update District column (G) for Chikkamagaluru rows (22-37)
$ws.Range("G22:G37").Value = "Chikkamagaluru (Chikmagalur)"

# Clear the stray empty inline string in F22 (Address column) so the cell no longer exists
$ws.Range("F22").ClearContents()
